$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# New row 43 - an extra export entry for Tabletop Simulator, following the
# same layout as the existing rows: p | Scheme | Mastermind | Villains |
# Henchmen | Heroes | Win | Scores | Close game | Special Rules | Notes
# (values are written in a specific column order so new shared strings are
# appended the same way the original workbook recorded them)
$row = 43
$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = "Secret Empire of Betrayal"
$ws.Cells.Item($row, 3).Value = "Mr. Sinister"
$ws.Cells.Item($row, 4).Value = "Marauders|Shadow-X"
$ws.Cells.Item($row, 5).Value = "Mandroid"
$ws.Cells.Item($row, 6).Value = "Deadpool (B)|Slapstick (DP)|Bullseye (V)|Nerkkod, Breaker of Oceans (FI)|Thing (FF)"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 11).Value = "Hit the final tactic turn before likely lethal twist."
$ws.Cells.Item($row, 10).Value = "War Machine as extra hero"
$ws.Cells.Item($row, 8).Value = "32|24"
$ws.Cells.Item($row, 9).Value = "yes"

$ws.Range("H44").Select()

# Scroll the view so D22 is the top-left visible cell (mirrors the
# topLeftCell="D22" recorded on the sheetView in the saved workbook).
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 4
